$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 13.16594766666667
$ws.Range("H2").Value = 39.497843
$ws.Range("I2").Value = 0.6940777873489595
$ws.Range("J2").Value = 0.6940777873489595
$ws.Range("M2").Value = 25.11140833333333
$ws.Range("N2").Value = 75.334225
$ws.Range("O2").Value = 0.7431105026796001
$ws.Range("P2").Value = 0.7431105026796001
$ws.Range("Q2").Value = 330.6154879529639
$ws.Range("R2").Value = 2975.539391576675
$ws.Range("S2").Value = 0.5157764934556299
$ws.Range("T2").Value = 0.5157764934556299

# Row 3
$ws.Range("G3").Value = 13.16594766666667
$ws.Range("H3").Value = 39.497843
$ws.Range("I3").Value = 0.6940777873489595
$ws.Range("J3").Value = 0.6940777873489595
$ws.Range("O3").Value = 0.1596166092346045
$ws.Range("P3").Value = 0.1596166092346045
$ws.Range("Q3").Value = 71.01463773854022
$ws.Range("R3").Value = 639.131739646862
$ws.Range("S3").Value = 0.1107863429616978
$ws.Range("T3").Value = 0.1107863429616978

# Row 4
$ws.Range("G4").Value = 13.16594766666667
$ws.Range("H4").Value = 39.497843
$ws.Range("I4").Value = 0.6940777873489595
$ws.Range("J4").Value = 0.6940777873489595
$ws.Range("N4").Value = 9.861222
$ws.Range("O4").Value = 0.09727288808579543
$ws.Range("P4").Value = 0.09727288808579541
$ws.Range("Q4").Value = 43.27744426046067
$ws.Range("R4").Value = 389.496998344146
$ws.Range("S4").Value = 0.06751495093163186
$ws.Range("T4").Value = 0.06751495093163185

# Row 5
$ws.Range("I5").Value = 0.1706596770095176
$ws.Range("J5").Value = 0.1706596770095176
$ws.Range("M5").Value = 25.11140833333333
$ws.Range("N5").Value = 75.334225
$ws.Range("O5").Value = 0.7431105026796001
$ws.Range("P5").Value = 0.7431105026796001
$ws.Range("Q5").Value = 81.29165551299999
$ws.Range("R5").Value = 731.624899617
$ws.Range("S5").Value = 0.1268189983696808
$ws.Range("T5").Value = 0.1268189983696808

# Row 6
$ws.Range("I6").Value = 0.1706596770095176
$ws.Range("J6").Value = 0.1706596770095176
$ws.Range("O6").Value = 0.1596166092346045
$ws.Range("P6").Value = 0.1596166092346045
$ws.Range("S6").Value = 0.02724011897733199
$ws.Range("T6").Value = 0.02724011897733199

# Row 7
$ws.Range("I7").Value = 0.1706596770095176
$ws.Range("J7").Value = 0.1706596770095176
$ws.Range("N7").Value = 9.861222
$ws.Range("O7").Value = 0.09727288808579543
$ws.Range("P7").Value = 0.09727288808579541
$ws.Range("R7").Value = 95.76942692183999
$ws.Range("S7").Value = 0.0166005596625048
$ws.Range("T7").Value = 0.0166005596625048

# Row 8
$ws.Range("H8").Value = 7.697376999999999
$ws.Range("I8").Value = 0.1352625356415228
$ws.Range("J8").Value = 0.1352625356415228
$ws.Range("M8").Value = 25.11140833333333
$ws.Range("N8").Value = 75.334225
$ws.Range("O8").Value = 0.7431105026796001
$ws.Range("P8").Value = 0.7431105026796001
$ws.Range("Q8").Value = 64.43065898086944
$ws.Range("R8").Value = 579.875930827825
$ws.Range("S8").Value = 0.1005150108542893
$ws.Range("T8").Value = 0.1005150108542893

# Row 9
$ws.Range("H9").Value = 7.697376999999999
$ws.Range("I9").Value = 0.1352625356415228
$ws.Range("J9").Value = 0.1352625356415228
$ws.Range("O9").Value = 0.1596166092346045
$ws.Range("P9").Value = 0.1596166092346045
$ws.Range("S9").Value = 0.0215901472955747
$ws.Range("T9").Value = 0.02159014729557471

# Row 10
$ws.Range("H10").Value = 7.697376999999999
$ws.Range("I10").Value = 0.1352625356415228
$ws.Range("J10").Value = 0.1352625356415228
$ws.Range("N10").Value = 9.861222
$ws.Range("O10").Value = 0.09727288808579543
$ws.Range("P10").Value = 0.09727288808579541
$ws.Range("Q10").Value = 8.433949268299333
$ws.Range("R10").Value = 75.90554341469399
$ws.Range("S10").Value = 0.01315737749165876
$ws.Range("T10").Value = 0.01315737749165876
